$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.703.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.338.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.337.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  +7.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.60%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.754.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.686.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("E17").Value = "  +2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.337.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.65%  "
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E26").Value = "  +6.00%  "
$ws.Range("E27").Value = "  +4.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.63%  "
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.931"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("E38").Value = "  +5.05%  "
$ws.Range("E39").Value = "  +7.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "274.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0931"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0502"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0216"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.71%  "
